$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "nerven"
$ws.Range("B3").Value = "dog/dog015.jpg"
$ws.Range("C3").Value = "dog"
$ws.Range("A4").Value = "rasen"
$ws.Range("B4").Value = "house/house004.jpg"
$ws.Range("C4").Value = "house"
$ws.Range("A6").Value = "geben"
$ws.Range("B6").Value = "dog/dog012.jpg"
$ws.Range("C6").Value = "dog"
$ws.Range("A7").Value = "fragen"
$ws.Range("B7").Value = "house/house008.jpg"
$ws.Range("C7").Value = "house"
$ws.Range("A9").Value = "stürmen"
$ws.Range("B9").Value = "dog/dog019.jpg"
$ws.Range("A10").Value = "küssen"
$ws.Range("B10").Value = "dog/dog020.jpg"
$ws.Range("A12").Value = "lächeln"
$ws.Range("B12").Value = "house/house012.jpg"
$ws.Range("A13").Value = "mühen"
$ws.Range("B13").Value = "house/house010.jpg"
$ws.Range("A15").Value = "tollen"
$ws.Range("B15").Value = "dog/dog023.jpg"
$ws.Range("C15").Value = "dog"
$ws.Range("A16").Value = "holen"
$ws.Range("B16").Value = "house/house026.jpg"
$ws.Range("C16").Value = "house"
$ws.Range("A18").Value = "machen"
$ws.Range("B18").Value = "dog/dog009.jpg"
$ws.Range("A19").Value = "zögern"
$ws.Range("B19").Value = "house/house009.jpg"
$ws.Range("C19").Value = "house"
$ws.Range("A21").Value = "spenden"
$ws.Range("B21").Value = "house/house021.jpg"
$ws.Range("A22").Value = "tragen"
$ws.Range("B22").Value = "dog/dog018.jpg"
$ws.Range("A24").Value = "schleppen"
$ws.Range("A25").Value = "süßen"
$ws.Range("B25").Value = "dog/dog022.jpg"
$ws.Range("C25").Value = "dog"
$ws.Range("A27").Value = "testen"
$ws.Range("B27").Value = "house/house013.jpg"
$ws.Range("A28").Value = "stillen"
$ws.Range("B28").Value = "dog/dog003.jpg"
$ws.Range("C28").Value = "dog"
$ws.Range("A30").Value = "bauen"
$ws.Range("B30").Value = "dog/dog028.jpg"
$ws.Range("C30").Value = "dog"
$ws.Range("A31").Value = "stärken"
$ws.Range("B31").Value = "dog/dog025.jpg"
$ws.Range("C31").Value = "dog"
$ws.Range("A33").Value = "binden"
$ws.Range("B33").Value = "dog/dog026.jpg"
$ws.Range("C33").Value = "dog"
$ws.Range("A34").Value = "quälen"
$ws.Range("B34").Value = "house/house031.jpg"
$ws.Range("C34").Value = "house"
$ws.Range("A36").Value = "leuchten"
$ws.Range("B36").Value = "house/house023.jpg"
$ws.Range("A37").Value = "schützen"
$ws.Range("B37").Value = "house/house019.jpg"
$ws.Range("A39").Value = "tauschen"
$ws.Range("B39").Value = "house/house006.jpg"
$ws.Range("C39").Value = "house"
$ws.Range("A40").Value = "kehren"
$ws.Range("B40").Value = "dog/dog017.jpg"
$ws.Range("C40").Value = "dog"
$ws.Range("A42").Value = "leugnen"
$ws.Range("B42").Value = "house/house030.jpg"
$ws.Range("C42").Value = "house"
$ws.Range("A43").Value = "packen"
$ws.Range("B43").Value = "house/house029.jpg"
$ws.Range("C43").Value = "house"
$ws.Range("A45").Value = "posten"
$ws.Range("B45").Value = "house/house025.jpg"
$ws.Range("C45").Value = "house"
$ws.Range("A46").Value = "steuern"
$ws.Range("B46").Value = "house/house018.jpg"
$ws.Range("C46").Value = "house"
$ws.Range("A48").Value = "kleben"
$ws.Range("B48").Value = "dog/dog016.jpg"
$ws.Range("A49").Value = "gelten"
$ws.Range("B49").Value = "dog/dog000.jpg"
$ws.Range("C49").Value = "dog"
